$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Re-import the "relasi" quiz data: 3 question rows, each followed by five
# A-E answer options (text in col B/D/F/H/J, a score in col C/E/G/I/K).
# ---------------------------------------------------------------------------

# Row 1 -----------------------------------------------------------------
$ws.Range("A1").Value = "TANGKAI : KELOPAK : BUNGA="
$ws.Range("B1").Value = "A. Tubuh : tangan : kepala"
$ws.Range("C1").Value = 10
$ws.Range("D1").Value = "B. Tanah : laut : air"
$ws.Range("E1").Value = 30
$ws.Range("F1").Value = "C. Tahun : bulan : hari"
$ws.Range("G1").Value = 20
$ws.Range("H1").Value = "D. Pelepah : tangkai : daun"
$ws.Range("I1").Value = 50
$ws.Range("J1").Value = "E. Langit : tanah : magma"
$ws.Range("K1").Value = 40

# Row 2 -----------------------------------------------------------------
$ws.Range("A2").Value = "A, B, C, F, E, D, G, H, I, L, …, …"
$ws.Range("B2").Value = "A. K dan J"
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = "B. J dan K"
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = "C. M dan N"
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = "D. N dan M"
$ws.Range("I2").Value = 40
$ws.Range("J2").Value = "E. I dan H"
$ws.Range("K2").Value = 30

# Row 3 -----------------------------------------------------------------
$ws.Range("A3").Value = "Perbandingan kelereng Egi dan Legi adalah 3 : 2. Jika selisih kelereng mereka 8, jumlah kelereng Egi dan Legi adalah…"
$ws.Range("B3").Value = "A. 40"
$ws.Range("C3").Value = 50
$ws.Range("D3").Value = "B. 32"
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = "C. 24"
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = "D. 16"
$ws.Range("I3").Value = 10
$ws.Range("J3").Value = "E. 25"
$ws.Range("K3").Value = 40

# ---------------------------------------------------------------------------
# Formatting: clear the old "centered" look from the question/number cells,
# and left-align the answer-option label cells (col B/D/F/H/J).
# ---------------------------------------------------------------------------
$ws.Range("A1:K3").Style = "Normal"
$leftAlignCells = "B1","D1","F1","H1","J1","B2","D2","F2","H2","J2","B3"
foreach ($addr in $leftAlignCells) {
    $ws.Range($addr).HorizontalAlignment = -4131
}

# Column widths (character units; chosen so the saved <col width> rounds to
# the same value Excel stored: 38.285.., 23.71.., 20.43.., 23.43.., 27.71..,
# 27.14..) ---------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 37.5
$ws.Columns("B:B").ColumnWidth = 22.833333333333332
$ws.Columns("D:D").ColumnWidth = 19.666666666666668
$ws.Columns("F:F").ColumnWidth = 22.666666666666668
$ws.Columns("H:H").ColumnWidth = 26.833333333333332
$ws.Columns("J:J").ColumnWidth = 26.333333333333332

# Selection + page setup, matching the re-saved workbook ---------------------
$ws.Columns("B:B").Select()
$ws.PageSetup.Orientation = 1
